# Data update using gitrun.py
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E15").Value = 97

$ws.Range("E18").Value = 49

$ws.Range("E36").Value = 46

$ws.Range("E41").Value = 17

$ws.Range("E43").Value = 12

$ws.Range("E46").Value = 13

$ws.Range("E48").Value = 11

$ws.Range("E49").Value = 32
$ws.Range("F49").Value = 15
$ws.Range("H49").Value = 15

$ws.Range("E50").Value = 9

$ws.Range("E64").Value = 21
$ws.Range("F64").Value = 11
$ws.Range("H64").Value = 11

$ws.Range("E65").Value = 16

$ws.Range("E67").Value = 19

$ws.Range("E71").Value = 12

$ws.Range("E74").Value = 9

$ws.Range("E77").Value = 22

$ws.Range("E83").Value = 5
$ws.Range("F83").Value = 1
$ws.Range("H83").Value = 1

$ws.Range("E88").Value = 7
